$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / tab-bar ratio -------------------------------------------------
# The original file's <workbookView tabRatio="993"/> grows to tabRatio="996".
# Set it through the window object (best-effort; some hosts treat TabRatio as
# a 0..1 fraction of the window width, others as the raw 0..1000 value, so
# try the fraction form first and the raw form as a fallback).
try {
    $excel.ActiveWindow.TabRatio = 0.996
} catch {
}
try {
    $excel.ActiveWindow.TabRatio = 996
} catch {
}

# --- Selection ---------------------------------------------------------------
# <selection activeCell="F62" .../> -> <selection activeCell="A59" .../>
$ws.Activate() | Out-Null
$ws.Range("A59").Select() | Out-Null

# --- Column widths -------------------------------------------------------
# Every data column got a little wider (same proportional bump LibreOffice
# applies when it recomputes "optimal" width). Re-apply the new widths
# (character units, i.e. the classic VBA `ColumnWidth`) column by column so
# each one lands on its new target width:
#   A/B : 19.5357142857143 -> 20.0255102040816
#   C   : 13.3571428571429 -> 13.7397959183673
#   D   : 23.1632653061224 -> 23.7602040816327
#   E   : 21.4030612244898 -> 21.984693877551
#   F   : 19.5357142857143 -> 20.0255102040816
#   G   : 33.5765306122449 -> 34.4591836734694
#   H.. : 19.5357142857143 -> 20.0255102040816
$ws.Columns.Item(1).ColumnWidth = 19.192176870748266
$ws.Columns.Item(2).ColumnWidth = 19.192176870748266
$ws.Columns.Item(3).ColumnWidth = 12.906462585033966
$ws.Columns.Item(4).ColumnWidth = 22.92687074829937
$ws.Columns.Item(5).ColumnWidth = 21.151360544217667
$ws.Columns.Item(6).ColumnWidth = 19.192176870748266
$ws.Columns.Item(7).ColumnWidth = 33.62585034013606

# Columns H:K are still inside the used range (A1:K59), so widen them to
# match the new default width used for the rest of the sheet.
$ws.Columns.Item(8).ColumnWidth = 19.192176870748266
$ws.Columns.Item(9).ColumnWidth = 19.192176870748266
$ws.Columns.Item(10).ColumnWidth = 19.192176870748266
$ws.Columns.Item(11).ColumnWidth = 19.192176870748266
